# QA and update files to console
# Applies the shared-string / style / layout changes described by the
# commit "QA and update files to console" to the computeList.xlsx
# workbook (sheet "ch").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C header text updates -------------------------------------
# Row 2, column C: "Intranet Address" -> "Private Address" (now shown in
# red to flag the re-translation for QA).
$ws.Range("C2").Value = "Private Address"
$ws.Range("C2").Font.Color = 255

# Row 4, column C: "No Virtual Machine Data" -> "No Virtual Machine data"
# (lower-case "data"), also flagged in red.
$ws.Range("C4").Value = "No Virtual Machine data"
$ws.Range("C4").Font.Color = 255

# --- Sheet view / selection --------------------------------------------
# The saved selection moves from D6 to I4.
$ws.Range("I4").Select() | Out-Null

# --- Page setup ----------------------------------------------------------
# A page setup block (portrait orientation) is now present for the sheet.
$ws.PageSetup.Orientation = 1
